$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Total_Risk2"
$ws.Range("B2").Value = "Test Total Risk with VARDEF=N"
